# The sheet holds one weekly price record per row (Cilantro, Agrícola del
# Norte S.A. de Arica). This commit adds one new weekly record, inserted
# as row 10 — every existing record from row 10 downward shifts down by
# one row (row 114 -> row 115), and the new dimension becomes A1:R115.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 10; rows 10..114 shift down to 11..115.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the new weekly record.
$ws.Range("A10").Value2 = 1
$ws.Range("B10").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value2 = "Arica y Parinacota"
$ws.Range("D10").Value2 = 45051
$ws.Range("D10").NumberFormat = $ws.Range("D11").NumberFormat
$ws.Range("E10").Value2 = 15
$ws.Range("F10").Value2 = 100112040
$ws.Range("G10").Value2 = "Cilantro"
$ws.Range("H10").Value2 = "Sin especificar"
$ws.Range("I10").Value2 = "Primera"
$ws.Range("J10").Value2 = 480
$ws.Range("K10").Value2 = 2800
$ws.Range("L10").Value2 = 3000
$ws.Range("M10").Value2 = 2917
$ws.Range("N10").Value2 = "$/atado 1,5 a 2 kilos"
$ws.Range("O10").Value2 = "Región de Arica y Parinacota"
$ws.Range("P10").Value2 = 1458
$ws.Range("Q10").Value2 = 2
$ws.Range("R10").Value2 = "Hortaliza"
